# Actualización automática 2025-08-19 13:10:10
$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# "VENTAS POR GRUPO": PORCELANATO sales (agosto) for RIOS CARRION ANGEL BENIGNO / CONZA VEGA FRANCO BLADYMIR
$wsGrupo.Range("M9").Value = 2161.81

# "VENTA MENSUAL": agosto column for the same client, plus the column total
$wsMensual.Range("F9").Value = 2161.81
$wsMensual.Range("F24").Value = 3219.08

# "CUMPLIMIENTO MENSUAL": PORCELANATO row (16) and TOTAL row (19)
$wsCumpl.Range("D16").Value = 3219.08
$wsCumpl.Range("E16").Value = 32837.62
$wsCumpl.Range("F16").Value = 0.08927827560481132

$wsCumpl.Range("D19").Value = 3219.08
$wsCumpl.Range("E19").Value = 51804.08386304603
$wsCumpl.Range("F19").Value = 0.05850408762412076

# Widen column F (CUMPLIMIENTO) from 24 to 25 characters
$wsCumpl.Columns.Item(6).ColumnWidth = 24.17
